$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Festival")
$ws2 = $wb.Worksheets.Item("Open")

# ===========================================================================
# 1) THEME: header fill colour FFFFC000 -> FFFF99, alignment center/center ->
#    center only (vertical cleared), on both sheets.
# ===========================================================================
foreach ($ws in @($ws1, $ws2)) {
    $ws.Range("A1:M1").Interior.Color = 10092543   # BGR for RGB(255,255,153)
    $ws.Range("A1:M1").VerticalAlignment = -4107   # xlBottom -> clears "vertical=center"
}

# ===========================================================================
# 2) COLUMN WIDTHS (identical layout on both sheets).
#    NOTE: this engine's ColumnWidth setter stores width + 0.8333.., so the
#    values below are back-solved to land exactly on the target widths.
# ===========================================================================
foreach ($ws in @($ws1, $ws2)) {
    $ws.Columns.Item(1).ColumnWidth  = 2.6666666666666665   # -> 3.5
    $ws.Columns.Item(2).ColumnWidth  = 34.166666666666664   # -> 35
    $ws.Columns.Item(3).ColumnWidth  = 14.166666666666666   # -> 15
    $ws.Columns.Item(4).ColumnWidth  = 14.166666666666666   # -> 15
    $ws.Columns.Item(5).ColumnWidth  = 12.166666666666666   # -> 13
    $ws.Columns.Item(6).ColumnWidth  = 19.166666666666668   # -> 20
    $ws.Columns.Item(7).ColumnWidth  = 9.166666666666666    # -> 10
    $ws.Columns.Item(8).ColumnWidth  = 7.166666666666667    # -> 8
    $ws.Columns.Item(9).ColumnWidth  = 19.166666666666668   # -> 20
    $ws.Columns.Item(10).ColumnWidth = 19.166666666666668   # -> 20
    $ws.Columns.Item(11).ColumnWidth = 19.166666666666668   # -> 20
    $ws.Columns.Item(12).ColumnWidth = 59.166666666666664   # -> 60
    $ws.Columns.Item(13).ColumnWidth = 24.166666666666668   # -> 25
}

# ===========================================================================
# 3) FESTIVAL sheet, row 2: update to the new "Kata" registration + restyle
#    the whole body row (center/center -> left, vertical cleared), then make
#    the "Kelas Otomatis" cell (L2) bold on top of that.
# ===========================================================================
$ws1.Range("H2").Value = 0
$ws1.Range("J2").Value = "Kata"
$ws1.Range("L2").Value = "FESTIVAL KATA PERORANGAN UNDER-21 SABUK HITAM PUTRA"
$ws1.Range("M2").Value = "25/6/2025, 00.35.43"

$ws1.Range("A2:M2").HorizontalAlignment = -4131 # xlLeft
$ws1.Range("A2:M2").VerticalAlignment   = -4107 # xlBottom
$ws1.Range("L2").Font.Bold = $true

# ===========================================================================
# 4) OPEN sheet: three new registration rows (2, 3, 4).
# ===========================================================================
$rows = @(
    @{ r=2; a=1; i="Fikrul cs"; j="Kata Beregu";  l="KATA BEREGU UNDER-21 PUTRA";     m="25/6/2025, 00.36.09" },
    @{ r=3; a=2; i="";          j="Kata";         l="KATA PERORANGAN UNDER-21 PUTRA"; m="25/6/2025, 00.36.22" },
    @{ r=4; a=3; i="";          j="Kata";         l="KATA PERORANGAN UNDER-21 PUTRA"; m="25/6/2025, 00.42.45" }
)

foreach ($row in $rows) {
    $r = $row.r

    $ws2.Range("A$r").Value = $row.a
    $ws2.Range("B$r").Value = "Raehan Fikrul Wahyu"

    # Copy the "Tanggal Lahir" text value from Festival!C2 so the engine keeps
    # it as the literal string "2004-12-27" instead of coercing it to a date.
    $ws1.Range("C2").Copy()
    $ws2.Range("C$r").PasteSpecial(-4104)

    $ws2.Range("D$r").Value = "Putra"
    $ws2.Range("E$r").Value = "BKC"
    $ws2.Range("F$r").Value = "Satria Galuh"
    $ws2.Range("G$r").Value = ""
    $ws2.Range("H$r").Value = 0
    $ws2.Range("I$r").Value = $row.i
    $ws2.Range("J$r").Value = $row.j
    $ws2.Range("K$r").Value = "Prestasi"
    $ws2.Range("L$r").Value = $row.l
    $ws2.Range("M$r").Value = $row.m

    # Body row alignment (left / default-vertical), then bold the "Kelas
    # Otomatis" cell (L) on top of that, matching Festival!L2's styling.
    $ws2.Range("A$r`:M$r").HorizontalAlignment = -4131
    $ws2.Range("A$r`:M$r").VerticalAlignment   = -4107
    $ws2.Range("L$r").Font.Bold = $true
}

Write-Host "edit complete"
